$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 655.2857
$ws.Range("I9").Value = 544.6667
$ws.Range("J9").Value = 738.25
$ws.Range("K9").Value = 544.6667
$ws.Range("L9").Value = 738.25
$ws.Range("M9").Value = -375.6667
$ws.Range("N9").Value = -1076.25
$ws.Range("H12").Value = 929.2
$ws.Range("I12").Value = 661.75
$ws.Range("J12").Value = 1999
$ws.Range("K12").Value = 661.75
$ws.Range("L12").Value = 1999
$ws.Range("M12").Value = -491.75
$ws.Range("N12").Value = -2339
$ws.Range("H38").Value = 247.6
$ws.Range("I38").Value = 164.22223
$ws.Range("J38").Value = 998
$ws.Range("K38").Value = 492.66669
$ws.Range("L38").Value = 2994
$ws.Range("M38").Value = -120.66669
$ws.Range("N38").Value = -3738
$ws.Range("H43").Value = 414527.6
$ws.Range("J43").Value = 460108.34
$ws.Range("L43").Value = 460108.34
$ws.Range("N43").Value = -460246.34
$ws.Range("H58").Value = 38466660
$ws.Range("I58").Value = 227.14285
$ws.Range("J58").Value = 83344170
$ws.Range("K58").Value = 681.4285500000001
$ws.Range("L58").Value = 250032510
$ws.Range("M58").Value = -531.4285500000001
$ws.Range("N58").Value = -250032810
$ws.Range("H92").Value = 425.9
$ws.Range("J92").Value = 303
$ws.Range("L92").Value = 303
$ws.Range("N92").Value = -2799
$ws.Range("H107").Value = 26786310
$ws.Range("I107").Value = 13889448
$ws.Range("J107").Value = 50000660
$ws.Range("K107").Value = 13889448
$ws.Range("L107").Value = 50000660
$ws.Range("M107").Value = -13887528
$ws.Range("N107").Value = -50004500
$ws.Range("H116").Value = 15631587
$ws.Range("I116").Value = 27781710
$ws.Range("J116").Value = 9999.286
$ws.Range("K116").Value = 27781710
$ws.Range("L116").Value = 9999.286
$ws.Range("M116").Value = -27778268
$ws.Range("N116").Value = -16883.286
$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178
$ws.Range("H118").Value = 483.16666
$ws.Range("I118").Value = 500
$ws.Range("K118").Value = 1500
$ws.Range("M118").Value = 157
$ws.Range("H132").Value = 1109.0244
$ws.Range("I132").Value = 1083.625
$ws.Range("K132").Value = 3250.875
$ws.Range("M132").Value = -720.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3034994.8
$ws.Range("I32").Value = 3229027
$ws.Range("J32").Value = 27495.5
$ws.Range("K32").Value = 3229027
$ws.Range("L32").Value = 27495.5
$ws.Range("M32").Value = -3228740
$ws.Range("N32").Value = -28069.5
$ws.Range("H74").Value = 42971.56
$ws.Range("I74").Value = 60545.59
$ws.Range("K74").Value = 60545.59
$ws.Range("M74").Value = -59671.59
$ws.Range("H77").Value = 42971.56
$ws.Range("I77").Value = 60545.59
$ws.Range("K77").Value = 302727.95
$ws.Range("M77").Value = -298359.95
$ws.Range("H102").Value = 2763
$ws.Range("I102").Value = 2605
$ws.Range("K102").Value = 2605
$ws.Range("M102").Value = -983
$ws.Range("H122").Value = 3959.7173
$ws.Range("I122").Value = 2966.2424
$ws.Range("J122").Value = 6481.615
$ws.Range("K122").Value = 8898.727200000001
$ws.Range("L122").Value = 19444.845
$ws.Range("M122").Value = -6448.727200000001
$ws.Range("N122").Value = -24344.845
$ws.Range("H132").Value = 6382.6587
$ws.Range("I132").Value = 3555.182
$ws.Range("J132").Value = 9656.579
$ws.Range("K132").Value = 10665.546
$ws.Range("L132").Value = 28969.737
$ws.Range("M132").Value = -8135.545999999998
$ws.Range("N132").Value = -34029.737

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5962662
$ws.Range("I107").Value = 8343727
$ws.Range("K107").Value = 8343727
$ws.Range("M107").Value = -8341807
$ws.Range("H135").Value = 79992.5
$ws.Range("J135").Value = 79992.5
$ws.Range("L135").Value = 79992.5
$ws.Range("N135").Value = -90132.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9796.096
$ws.Range("I58").Value = 3236.3333
$ws.Range("K58").Value = 3236.3333
$ws.Range("M58").Value = -3033.3333
$ws.Range("H107").Value = 1702.55
$ws.Range("J107").Value = 3219.5
$ws.Range("L107").Value = 3219.5
$ws.Range("N107").Value = -7059.5
$ws.Range("H136").Value = 9796.096
$ws.Range("I136").Value = 3236.3333
$ws.Range("K136").Value = 9708.999899999999
$ws.Range("M136").Value = -7158.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2172.3635
$ws.Range("I131").Value = 1099.1333
$ws.Range("K131").Value = 3297.3999
$ws.Range("M131").Value = 1742.6001
$ws.Range("H138").Value = 66417.5
$ws.Range("I138").Value = 81206.38
$ws.Range("K138").Value = 243619.14
$ws.Range("M138").Value = -238479.14

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 737454.4399999999
$ws.Range("I70").Value = 8000000
$ws.Range("J70").Value = 11199.9
$ws.Range("K70").Value = 8000000
$ws.Range("L70").Value = 11199.9
$ws.Range("M70").Value = -7999730
$ws.Range("N70").Value = -11739.9
$ws.Range("H73").Value = 737454.4399999999
$ws.Range("I73").Value = 8000000
$ws.Range("J73").Value = 11199.9
$ws.Range("K73").Value = 8000000
$ws.Range("L73").Value = 11199.9
$ws.Range("M73").Value = -7999064
$ws.Range("N73").Value = -13071.9
$ws.Range("H113").Value = 351234
$ws.Range("I113").Value = 911872.2
$ws.Range("J113").Value = 8621.777
$ws.Range("K113").Value = 911872.2
$ws.Range("L113").Value = 8621.777
$ws.Range("M113").Value = -909702.2
$ws.Range("N113").Value = -12961.777
$ws.Range("H132").Value = 6542
$ws.Range("I132").Value = 2807.4285
$ws.Range("K132").Value = 8422.2855
$ws.Range("M132").Value = -5892.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3392.8572
$ws.Range("J22").Value = 3933.3333
$ws.Range("L22").Value = 3933.3333
$ws.Range("N22").Value = -4523.3333
$ws.Range("H27").Value = 3392.8572
$ws.Range("J27").Value = 3933.3333
$ws.Range("L27").Value = 3933.3333
$ws.Range("N27").Value = -4147.3333
$ws.Range("H40").Value = 18523312
$ws.Range("I40").Value = 25003696
$ws.Range("J40").Value = 7929.2856
$ws.Range("K40").Value = 25003696
$ws.Range("L40").Value = 7929.2856
$ws.Range("M40").Value = -25003560
$ws.Range("N40").Value = -8201.285599999999
$ws.Range("H132").Value = 8781.351000000001
$ws.Range("I132").Value = 2756.9333
$ws.Range("K132").Value = 8270.7999
$ws.Range("M132").Value = -5740.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10753417
$ws.Range("I107").Value = 282.2353
$ws.Range("K107").Value = 846.7058999999999
$ws.Range("M107").Value = 1073.2941
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H136").Value = 316473.94
$ws.Range("J136").Value = 674180.2
$ws.Range("L136").Value = 2022540.6
